# Add ZL plans for sprint2
#
# On the "Sprint2" sheet, the US03 ("Birth before death") story block and
# the US05 ("Marriage before death") story block are replaced with two new
# ZL-owned stories: US22 ("Unique IDs") and US29 ("List deceased"), each
# with two sub-tasks. Both new story blocks are one row shorter than the
# blocks they replace (2 tasks instead of 3), so the sheet shrinks by two
# rows overall and every later row shifts up accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")

# Drop one task row from each of the two blocks being replaced so the
# remaining rows collapse into the new (shorter) layout; everything below
# shifts up automatically.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(17).Delete()

# Normalize row heights across the rewritten block (clears the taller
# wrapped-text rows left over from the old multi-line task descriptions).
$ws.Range("A10:I18").Rows.AutoFit()

# --- US22: Unique IDs (replaces the old US03 block) ---
$ws.Range("A10").Value = "US22"
$ws.Range("B10").Value = "Unique IDs"
$ws.Range("D10").Value = "Coding"
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = 30

$ws.Range("A11").Value = "T22.01"
$ws.Range("B11").Value = "parse all IDs"

$ws.Range("A12").Value = "T22.02"
$ws.Range("B12").Value = "find duplicates"

# --- US29: List deceased (replaces the old US05 block) ---
$ws.Range("A14").Value = "US29"
$ws.Range("B14").Value = "List deceased"
$ws.Range("D14").Value = "Coding"
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 20

$ws.Range("A15").Value = "T29.01"
$ws.Range("B15").Value = "determine sort standard"

$ws.Range("A16").Value = "T29.02"
$ws.Range("B16").Value = "sort"

# The Sprint2 tab becomes the active tab/sheet, with G13 selected (Sprint1
# loses tabSelected automatically since only one sheet can carry it).
$ws.Activate()
$ws.Range("G13").Select()
